$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new shared-string values in the exact order needed so the
# shared-strings table grows with the same index assignment as the
# authored change (355 "Effekseer Project" .. 363 "StorageTarget").
$ws.Range("B144").Value = "Effekseer Project"
$ws.Range("B142").Value = "Setting save destination"
$ws.Range("C142").Value = "設定保存先"
$ws.Range("B143").Value = "Effekseer Application"
$ws.Range("C143").Value = "Effekseerアプリケーション"
$ws.Range("C144").Value = "プロジェクト"
$ws.Range("A143").Value = "StorageGlobal"
$ws.Range("A144").Value = "StorageLocal"
$ws.Range("A142").Value = "StorageTarget"

# Apply the same wrap-text / vertically-centered style used throughout the
# sheet ("s=1") to every new cell except B144, which the author left
# unstyled (default style), matching the rest of the table's convention.
$styled = $ws.Range("A142:C143")
$styled.WrapText = $true
$styled.VerticalAlignment = -4108

$ws.Range("A144").WrapText = $true
$ws.Range("A144").VerticalAlignment = -4108
$ws.Range("C144").WrapText = $true
$ws.Range("C144").VerticalAlignment = -4108

# Columns got very slightly re-measured by Excel when the new rows were
# added; reproduce the closest attainable widths (the host's ColumnWidth
# setter quantizes to 1/7-character steps, so these land on the nearest
# representable value to the authored 34.42578125 / 80.28515625).
$ws.Columns.Item(1).ColumnWidth = 33.65
$ws.Columns.Item(2).ColumnWidth = 79.55

# Move the selection/scroll position to reflect where the new rows were
# authored.
$null = $ws.Range("A142").Select()
